# Add DBV control settings to the epanet model's flow device database.
#
# The physical PRV/DBV flow meters at Snowden Road and New Station Way
# (and a handful of related FM/inlet rows that log against the same
# devices) were re-pointed to their DBV-control logger IDs. Some of the
# new IDs are alphanumeric ("<number>_<suffix>") so Excel stores them as
# text, others are still plain numeric device IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  - Stoke Lane PRV / inlet_2296 logger -> new DBV-style asset id (text)
$ws.Range("J2").Value = "620104099_2"

# Row 3  - Woodland Way PRV / inlet_2005 logger -> new DBV-style asset id (text)
$ws.Range("J3").Value = "644926906_2"

# Row 4  - Lodge Causeway PRV -> new numeric device id
$ws.Range("J4").Value = 579829886

# Row 5  - Snowden Road DBV -> new DBV-style asset id (text)
$ws.Range("J5").Value = "623990497_2"

# Row 6  - New Station Way DBV -> new numeric device id
$ws.Range("J6").Value = 623990697

# Row 7  - inlet_2296 flow logger -> new DBV-style asset id (text)
$ws.Range("J7").Value = "612058339_1"

# Row 8  - inlet_2005 flow logger -> new DBV-style asset id (text)
$ws.Range("J8").Value = "644926906_2"

# Row 9  - wastemeter_2297 -> new numeric device id
$ws.Range("J9").Value = 733663782

# Row 14 - wastemeter_2302 -> new numeric device id
$ws.Range("J14").Value = 579829886

# Row 16 - wastemeter_2306 (inactive row) -> asset id cleared entirely
$ws.Range("J16").Clear()

# Row 17 - wastemeter_2308 -> new DBV-style asset id (text)
$ws.Range("J17").Value = "623990497_2"

# Row 18 - wastemeter_2307 -> new numeric device id
$ws.Range("J18").Value = 623990697

# Move/restore the active selection as left by the editor
$ws.Range("M26").Select()
